# Natmi following Dr Hou advice
# Rewrites the LR-pairs data block (rows 2-17) with a new 4x4 sending/target
# cluster grid (ECs added as a sending cluster) and refreshed statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Adam23","Itgb3","ECs",2,0.6666666666666666,0.269164,0.8074920000000001,0.01089095165781685,0.01089095165781686,3,1,3.778439,11.335317,0.4252971528324392,0.4252971528324392,1.017019754996,9.153177794964002,0.004631890731705242,0.004631890731705243),
    @("ECs","Adam23","Itgb3","FAPs",2,0.6666666666666666,0.269164,0.8074920000000001,0.01089095165781685,0.01089095165781686,3,1,4.333403333333333,13.00021,0.4877633593505858,0.4877633593505858,1.166396174813333,10.49756557332,0.00531220716714158,0.005312207167141582),
    @("ECs","Adam23","Itgb3","M2",2,0.6666666666666666,0.269164,0.8074920000000001,0.01089095165781685,0.01089095165781686,3,1,0.2909853333333334,0.8729560000000001,0.03275300561492853,0.03275300561492853,0.07832277626133335,0.7049049863520002,0.0003567114008003907,0.0003567114008003907),
    @("ECs","Adam23","Itgb3","sCs",2,0.6666666666666666,0.269164,0.8074920000000001,0.01089095165781685,0.01089095165781686,3,1,0.4814053333333333,1.444216,0.0541864822020464,0.05418648220204641,0.1295769851413333,1.166192866272,0.0005901423581696407,0.0005901423581696409),
    @("FAPs","Adam23","Itgb3","ECs",3,1,20.07911633333334,60.237349,0.8124440315867432,0.8124440315867433,3,1,3.778439,11.335317,0.4252971528324392,0.4252971528324392,75.86771623940368,682.8094461546331,0.3455301334695502,0.3455301334695502),
    @("FAPs","Adam23","Itgb3","FAPs",3,1,20.07911633333334,60.237349,0.8124440315867432,0.8124440315867433,3,1,4.333403333333333,13.00021,0.4877633593505858,0.4877633593505858,87.01090964925444,783.09818684329,0.3962804301310833,0.3962804301310834),
    @("FAPs","Adam23","Itgb3","M2",3,1,20.07911633333334,60.237349,0.8124440315867432,0.8124440315867433,3,1,0.2909853333333334,0.8729560000000001,0.03275300561492853,0.03275300561492853,5.842728359293779,52.584555233644,0.02660998392837577,0.02660998392837578),
    @("FAPs","Adam23","Itgb3","sCs",3,1,20.07911633333334,60.237349,0.8124440315867432,0.8124440315867433,3,1,0.4814053333333333,1.444216,0.0541864822020464,0.05418648220204641,9.666193691487111,86.995743223384,0.04402348405773388,0.0440234840577339),
    @("M2","Adam23","Itgb3","ECs",3,1,0.174593,0.523779,0.007064406543197522,0.007064406543197523,3,1,3.778439,11.335317,0.4252971528324392,0.4252971528324392,0.6596890003270001,5.937201002943,0.00300447198927276,0.003004471989272761),
    @("M2","Adam23","Itgb3","FAPs",3,1,0.174593,0.523779,0.007064406543197522,0.007064406543197523,3,1,4.333403333333333,13.00021,0.4877633593505858,0.4877633593505858,0.7565818881766666,6.80923699359,0.003445758667328283,0.003445758667328283),
    @("M2","Adam23","Itgb3","M2",3,1,0.174593,0.523779,0.007064406543197522,0.007064406543197523,3,1,0.2909853333333334,0.8729560000000001,0.03275300561492853,0.03275300561492853,0.05080400230266667,0.457236020724,0.0002313805471754863,0.0002313805471754863),
    @("M2","Adam23","Itgb3","sCs",3,1,0.174593,0.523779,0.007064406543197522,0.007064406543197523,3,1,0.4814053333333333,1.444216,0.0541864822020464,0.05418648220204641,0.08405000136266666,0.756450012264,0.0003827953394209927,0.0003827953394209928),
    @("sCs","Adam23","Itgb3","ECs",3,1,4.191587666666666,12.574763,0.1696006102122423,0.1696006102122424,3,1,3.778439,11.335317,0.4252971528324392,0.4252971528324392,15.83765831165233,142.538924804871,0.07213065664191098,0.07213065664191101),
    @("sCs","Adam23","Itgb3","FAPs",3,1,4.191587666666666,12.574763,0.1696006102122423,0.1696006102122424,3,1,4.333403333333333,13.00021,0.4877633593505858,0.4877633593505858,18.16383996669222,163.47455970023,0.08272496338503259,0.08272496338503263),
    @("sCs","Adam23","Itgb3","M2",3,1,4.191587666666666,12.574763,0.1696006102122423,0.1696006102122424,3,1,0.2909853333333334,0.8729560000000001,0.03275300561492853,0.03275300561492853,1.219690534380889,10.977214809428,0.005554929738576879,0.005554929738576881),
    @("sCs","Adam23","Itgb3","sCs",3,1,4.191587666666666,12.574763,0.1696006102122423,0.1696006102122424,3,1,0.4814053333333333,1.444216,0.0541864822020464,0.05418648220204641,2.017852657867555,18.160673920808,0.00919006044672188,0.009190060446721883)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
